$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 295
$ws.Range("I2").Value = 295
$ws.Range("K2").Value = 295
$ws.Range("M2").Value = -182

$ws.Range("H19").Value = 1222.2354
$ws.Range("I19").Value = 948.75
$ws.Range("K19").Value = 948.75
$ws.Range("M19").Value = -773.75

$ws.Range("H98").Value = 3457.1052
$ws.Range("I98").Value = 3087.353
$ws.Range("K98").Value = 3087.353
$ws.Range("M98").Value = -1589.353

$ws.Range("H106").Value = 6727.25
$ws.Range("I106").Value = 3948.5
$ws.Range("K106").Value = 3948.5
$ws.Range("M106").Value = -3317.5

$ws.Range("H113").Value = 23097.643
$ws.Range("I113").Value = 34507.445
$ws.Range("K113").Value = 34507.445
$ws.Range("M113").Value = -31253.445

$ws.Range("H122").Value = 3457.1052
$ws.Range("I122").Value = 3087.353
$ws.Range("K122").Value = 9262.059000000001
$ws.Range("M122").Value = -6812.059000000001

$ws.Range("H132").Value = 1316.7916
$ws.Range("I132").Value = 1079.439
$ws.Range("K132").Value = 3238.317
$ws.Range("M132").Value = -708.317

$ws.Range("H135").Value = 436.58066
$ws.Range("I135").Value = 436.58066
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3929.22594
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1394.22594
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 26918.615
$ws.Range("I137").Value = 1135.4348
$ws.Range("J137").Value = 63981.938
$ws.Range("K137").Value = 3406.3044
$ws.Range("L137").Value = 191945.814
$ws.Range("M137").Value = -856.3044
$ws.Range("N137").Value = -197045.814

$ws.Range("H138").Value = 1900.3889
$ws.Range("I138").Value = 1628.5084
$ws.Range("J138").Value = 2417.8386
$ws.Range("K138").Value = 4885.5252
$ws.Range("L138").Value = 7253.5158
$ws.Range("M138").Value = 254.4748
$ws.Range("N138").Value = -17533.5158

$ws.Range("H141").Value = 779539.4399999999
$ws.Range("I141").Value = 849362.5600000001
$ws.Range("K141").Value = 2548087.68
$ws.Range("M141").Value = -2542907.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3413.6736
$ws.Range("I32").Value = 2836.988
$ws.Range("K32").Value = 2836.988
$ws.Range("M32").Value = -2549.988

$ws.Range("H61").Value = 4471.95
$ws.Range("I61").Value = 1433.8462
$ws.Range("J61").Value = 10114.143
$ws.Range("K61").Value = 1433.8462
$ws.Range("L61").Value = 10114.143
$ws.Range("M61").Value = -1221.8462
$ws.Range("N61").Value = -10538.143

$ws.Range("H63").Value = 6699.8
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 6699.8
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H97").Value = 880.6
$ws.Range("I97").Value = 875.94446
$ws.Range("J97").Value = 922.5
$ws.Range("K97").Value = 875.94446
$ws.Range("L97").Value = 922.5
$ws.Range("M97").Value = -379.94446
$ws.Range("N97").Value = -1914.5

$ws.Range("H136").Value = 4471.95
$ws.Range("I136").Value = 1433.8462
$ws.Range("J136").Value = 10114.143
$ws.Range("K136").Value = 4301.5386
$ws.Range("L136").Value = 30342.429
$ws.Range("M136").Value = -1751.5386
$ws.Range("N136").Value = -35442.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H94").Value = 1622.6666
$ws.Range("I94").Value = 2287.6667
$ws.Range("J94").Value = 292.66666
$ws.Range("K94").Value = 2287.6667
$ws.Range("L94").Value = 292.66666
$ws.Range("M94").Value = -1836.6667
$ws.Range("N94").Value = -1194.66666

$ws.Range("H99").Value = 1517.3334
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1517.3334
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1517.3334
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -4513.3334

$ws.Range("H134").Value = 4463.657
$ws.Range("I134").Value = 4529.9355
$ws.Range("K134").Value = 13589.8065
$ws.Range("M134").Value = -11054.8065

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1901.8889
$ws.Range("I31").Value = 1388.8636
$ws.Range("K31").Value = 1388.8636
$ws.Range("M31").Value = -1093.8636

$ws.Range("H34").Value = 1901.8889
$ws.Range("I34").Value = 1388.8636
$ws.Range("K34").Value = 1388.8636
$ws.Range("M34").Value = -1186.8636

$ws.Range("H132").Value = 2008.9231
$ws.Range("I132").Value = 1360.7931
$ws.Range("J132").Value = 3888.5
$ws.Range("K132").Value = 4082.379300000001
$ws.Range("L132").Value = 11665.5
$ws.Range("M132").Value = -1552.379300000001
$ws.Range("N132").Value = -16725.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1250
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1250
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3750
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3932

$ws.Range("H58").Value = 2466.3333
$ws.Range("I58").Value = 2450
$ws.Range("J58").Value = 2499
$ws.Range("K58").Value = 7350
$ws.Range("L58").Value = 7497
$ws.Range("M58").Value = -7222
$ws.Range("N58").Value = -7753

$ws.Range("H121").Value = 789.8570999999999
$ws.Range("I121").Value = 815
$ws.Range("J121").Value = 779.8
$ws.Range("K121").Value = 2445
$ws.Range("L121").Value = 2339.4
$ws.Range("M121").Value = -1135
$ws.Range("N121").Value = -4959.4

$ws.Range("H122").Value = 1541.0435
$ws.Range("J122").Value = 1774.3334
$ws.Range("L122").Value = 15969.0006
$ws.Range("N122").Value = -20869.0006

$ws.Range("H131").Value = 19266022
$ws.Range("I131").Value = 83333950
$ws.Range("J131").Value = 45644.9
$ws.Range("K131").Value = 250001850
$ws.Range("L131").Value = 136934.7
$ws.Range("M131").Value = -249996810
$ws.Range("N131").Value = -147014.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1089.5
$ws.Range("I113").Value = 1027.7142
$ws.Range("J113").Value = 1176
$ws.Range("K113").Value = 1027.7142
$ws.Range("L113").Value = 1176
$ws.Range("M113").Value = 1142.2858
$ws.Range("N113").Value = -5516

$ws.Range("H126").Value = 1490560.5
$ws.Range("I126").Value = 3474365.5
$ws.Range("K126").Value = 10423096.5
$ws.Range("M126").Value = -10420626.5

$ws.Range("H132").Value = 1133811
$ws.Range("I132").Value = 1604553.4
$ws.Range("J132").Value = 4029.1
$ws.Range("K132").Value = 4813660.199999999
$ws.Range("L132").Value = 12087.3
$ws.Range("M132").Value = -4811130.199999999
$ws.Range("N132").Value = -17147.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2376.6667
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 2252
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 2252
$ws.Range("M22").Value = -2705
$ws.Range("N22").Value = -2842

$ws.Range("H27").Value = 2376.6667
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 2252
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 2252
$ws.Range("M27").Value = -2893
$ws.Range("N27").Value = -2466

$ws.Range("H46").Value = 1629.4286
$ws.Range("I46").Value = 922.75
$ws.Range("J46").Value = 2571.6667
$ws.Range("K46").Value = 922.75
$ws.Range("L46").Value = 2571.6667
$ws.Range("M46").Value = -734.75
$ws.Range("N46").Value = -2947.6667

$ws.Range("H82").Value = 2002
$ws.Range("I82").Value = 1475
$ws.Range("K82").Value = 1475
$ws.Range("M82").Value = -1114

$ws.Range("H85").Value = 2002
$ws.Range("I85").Value = 1475
$ws.Range("K85").Value = 1475
$ws.Range("M85").Value = -227

$ws.Range("H96").Value = 85000
$ws.Range("J96").Value = 85000
$ws.Range("L96").Value = 85000
$ws.Range("N96").Value = -90492

$ws.Range("H132").Value = 1971.5646
$ws.Range("I132").Value = 1380.7046
$ws.Range("J132").Value = 3415.889
$ws.Range("K132").Value = 4142.1138
$ws.Range("L132").Value = 10247.667
$ws.Range("M132").Value = -1612.1138
$ws.Range("N132").Value = -15307.667

$ws.Range("H136").Value = 1966.7656
$ws.Range("I136").Value = 1391.3729
$ws.Range("J136").Value = 8756.4
$ws.Range("K136").Value = 4174.1187
$ws.Range("L136").Value = 26269.2
$ws.Range("M136").Value = -1624.1187
$ws.Range("N136").Value = -31369.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6546.64
$ws.Range("I126").Value = 7486.294
$ws.Range("K126").Value = 22458.882
$ws.Range("M126").Value = -19988.882

$ws.Range("H132").Value = 1744.5476
$ws.Range("I132").Value = 1314.5714
$ws.Range("J132").Value = 2604.5
$ws.Range("K132").Value = 3943.7142
$ws.Range("L132").Value = 7813.5
$ws.Range("M132").Value = -1413.7142
$ws.Range("N132").Value = -12873.5
